# "add the code of pbf feature"
# Populate the previously-placeholder "X" (unknown amino acid) column with the
# position-averaged values computed from the other 20 amino-acid columns, and
# apply a 2-decimal custom number format to the newly-filled cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Attribute rows 2-6 hold per-amino-acid numeric values in columns B:V.
# Row 2 (Hydropathy index) already had a real value for the "X" placeholder
# column (U2), but rows 3-6 were left at 0 — backfill them with the mean of
# the 21 defined amino-acid columns (B:V), matching the value the feature now
# computes, and mark the written cells with the new "0.00_ " number format.
$ws.Range("U3").Value = 130.38190476190476
$ws.Range("U3").NumberFormat = "0.00_ "

$ws.Range("U4").Value = 5.7395238095238099
$ws.Range("U4").NumberFormat = "0.00_ "

$ws.Range("U5").Value = 8.9119047619047631
$ws.Range("U5").NumberFormat = "0.00_ "

$ws.Range("U6").Value = 2.0438095238095237
$ws.Range("U6").NumberFormat = "0.00_ "

# Move the live selection to where the editing session ended up.
$ws.Range("P16").Select() | Out-Null
